# Auto-generated edit script applying cryptos.xlsx price/volume refresh
# (GitHub Actions style data update: Sat Jul 13 03:34:16 UTC 2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}


Set-TextCell 'D2' '57.869.98'
Set-TextCell 'E2' '  +1.43%  '
Set-TextCell 'D3' '3.126.06'
Set-TextCell 'E3' '  +1.20%  '
Set-TextCell 'E4' '  -0.01%  '
Set-TextCell 'D5' '534.52'
Set-TextCell 'E5' '  +2.27%  '
Set-TextCell 'D6' '139.21'
Set-TextCell 'E6' '  +1.67%  '
Set-TextCell 'E7' '  -0.05%  '
Set-TextCell 'D8' '0.491'
Set-TextCell 'E8' '  +8.62%  '
Set-TextCell 'D9' '7.35'
Set-TextCell 'E9' '  -0.18%  '
Set-TextCell 'E10' '  +1.29%  '
Set-TextCell 'E11' '  +3.01%  '
Set-TextCell 'D12' '0.140'
Set-TextCell 'E12' '  +3.66%  '
Set-TextCell 'D13' '3.665.91'
Set-TextCell 'E13' '  +1.15%  '
Set-TextCell 'D14' '25.90'
Set-TextCell 'E14' '  +1.50%  '
Set-TextCell 'E15' '  +2.93%  '
Set-TextCell 'D16' '57.993.04'
Set-TextCell 'E16' '  +1.46%  '
Set-TextCell 'D17' '3.126.28'
Set-TextCell 'E17' '  +1.22%  '
Set-TextCell 'D18' '6.11'
Set-TextCell 'E18' '  +3.46%  '
Set-TextCell 'E19' '  +2.70%  '
Set-TextCell 'E20' '  +3.00%  '
Set-TextCell 'D21' '375.40'
Set-TextCell 'E23' '  -1.11%  '
Set-TextCell 'E24' '  +1.96%  '
Set-TextCell 'E25' '  +1.67%  '
Set-TextCell 'E26' '  -0.75%  '
Set-TextCell 'E27' '  -0.06%  '
Set-TextCell 'D28' '0.0₃0880'
Set-TextCell 'E28' '  -0.26%  '
Set-TextCell 'D29' '7.53'
Set-TextCell 'E29' '  +2.93%  '
Set-TextCell 'D30' '6.16'
Set-TextCell 'E30' '  +4.22%  '
Set-TextCell 'E31' '  +0.10%  '
Set-TextCell 'D32' '21.57'
Set-TextCell 'E32' '  +3.65%  '
Set-TextCell 'D33' '5.18'
Set-TextCell 'E33' '  +3.78%  '
Set-TextCell 'D35' '160.24'
Set-TextCell 'E35' '  +0.58%  '
Set-TextCell 'D36' '6.16'
Set-TextCell 'E37' '  +4.09%  '
Set-TextCell 'D38' '25.60'
Set-TextCell 'E38' '  -1.66%  '
Set-TextCell 'D39' '1.64'
Set-TextCell 'E39' '  +4.31%  '
Set-TextCell 'E40' '  +2.53%  '
Set-TextCell 'D41' '2.566.14'
Set-TextCell 'E41' '  +7.45%  '
Set-TextCell 'D42' '4.15'
Set-TextCell 'E42' '  +2.65%  '
Set-TextCell 'D43' '0.701'
Set-TextCell 'E43' '  +0.90%  '
Set-TextCell 'D44' '38.09'
Set-TextCell 'E44' '  +3.99%  '
Set-TextCell 'E45' '  +2.77%  '
Set-TextCell 'E46' '  +0.00%  '
Set-TextCell 'E47' '  +1.88%  '
Set-TextCell 'E48' '  +3.05%  '
Set-TextCell 'D49' '19.84'
Set-TextCell 'E49' '  +0.89%  '
Set-TextCell 'B50' 'Stellar'
Set-TextCell 'C50' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D50' '0.0934'
Set-TextCell 'E50' '  +4.50%  '
Set-TextCell 'B51' 'SuiNetwork'
Set-TextCell 'C51' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell 'D51' '0.749'
Set-TextCell 'E51' '  -1.79%  '
